$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column C
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 2
$ws.Range("C7").Value = 1

# Update the selected cell to A2
$ws.Range("A2").Select()
